$d = $word.ActiveDocument

# Update the header date text from "December 2022" to "August 2020".
# There are two occurrences inside the header (one in the drawing
# canvas content, one in the VML fallback), both need the same change.
foreach ($sec in $d.Sections) {
    $hdr = $sec.Headers.Item(1)  # wdHeaderFooterPrimary = 1
    if ($hdr.Exists) {
        $hdr.Range.Find.Execute("December 2022", $false, $false, $false, $false, $false, $true, 1, $false, "August 2020", 2) | Out-Null
    }
}
